# Apply the "request-issues" workbook update:
# - rename header C1 from "resolved" to "resolved_at_PH_end"
# - simplify the FALSE explanation strings in C8/C9
# - mark F8/F9 (resolved_at_PH_end) as TRUE instead of FALSE
# - move the active cell selection to D13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename
$ws.Range("C1").Value = "resolved_at_PH_end"

# Row 8 (request 938 / question_id 938, raised 43172)
$ws.Range("C8").Value = "FALSE, last month of data are dropped (maybe?)"
$ws.Range("F8").Value = $true

# Row 9 (request 1055)
$ws.Range("C9").Value = "FALSE, last month of data are dropped"
$ws.Range("F9").Value = $true

# Update the last active selection to match the saved file
$ws.Range("D13").Select()
